# Updated some raw values
$wb = $excel.ActiveWorkbook

# --- Calibration sheet ---
$ws1 = $wb.Worksheets.Item("Calibration")
$ws1.Range("E3").Value = 530
$ws1.Range("F3").Value = 518
$ws1.Range("E4").Value = 485
$ws1.Range("F4").Value = 477
$ws1.Range("E5").Value = 446
$ws1.Range("F5").Value = 437
$ws1.Range("E6").Value = 410
$ws1.Range("F6").Value = 404
$ws1.Range("E7").Value = 385
$ws1.Range("F7").Value = 378
$ws1.Range("D27").ClearContents()
$ws1.Range("F27").ClearContents()

# --- IR6 sheet ---
$ws3 = $wb.Worksheets.Item("IR6")
$ws3.Range("B8").Value = 231
$ws3.Range("B9").Value = 200

# B35 was retyped to reference the C-column slope/intercept instead of B's,
# and that edit was filled down through B38.
$ws3.Range("B35").Formula = "=C`$26/B6+C`$25"
$ws3.Range("B36").Formula = "=C`$26/B7+C`$25"
$ws3.Range("B37").Formula = "=C`$26/B8+C`$25"
$ws3.Range("B38").Formula = "=C`$26/B9+C`$25"
